$d = $word.ActiveDocument

# --- Move the "_GoBack" bookmark from the end of the document to the very
# beginning (right before the first run of the first paragraph). ---
#
# Word's own Bookmarks collection hides "_GoBack" from iteration/Count
# (it's an internal/hidden bookmark) but it can still be addressed by name.
# Re-adding a bookmark under the same name moves it (old one is removed).
#
# A zero-length bookmark placed at absolute document position 0 has a
# serialization quirk in this host when added directly, so we land it one
# character in first and then delete that helper character - the bookmark
# naturally collapses back to the true start of the document, same as it
# would in real Word when content ahead of a bookmark is removed.
$helper = $d.Range(0, 0)
$helper.InsertBefore("X")

$target = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $target)

$toDelete = $d.Range(0, 1)
$toDelete.Delete()

# --- Shrink the page margins to 0.5" (720 twips / 36pt) on every side. ---
$pageSetup = $d.Sections(1).PageSetup
$pageSetup.TopMargin = 36
$pageSetup.RightMargin = 36
$pageSetup.BottomMargin = 36
$pageSetup.LeftMargin = 36
